$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "AUS"
$ws.Range("C1").Value = "CHN"
$ws.Range("D1").Value = "DEU"
$ws.Range("E1").Value = "FRA"
$ws.Range("F1").Value = "USA"
